$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the two runs "THU Mar 07" + " 14:37:12 IST 2019"
# (first record in the log) into a single run with the combined text.
# -----------------------------------------------------------------
$found = $d.Content.Find.Execute("THU Mar 07 14:37:12 IST 2019", $false, $false, $false, $false, $false, $true, 1, $false, "THU Mar 07 14:37:12 IST 2019", 2)

# -----------------------------------------------------------------
# Change 2: append a brand new purchase record (FRI Mar 08 ...) right
# after the LAST "Amount balance" paragraph (the CHANDRA/BEET entry
# dated 08/03/2019), pushing the existing trailing blank paragraphs
# further down (they are left untouched).
# -----------------------------------------------------------------

# Locate the last paragraph whose text contains "Amount balance".
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Amount balance*") {
        $targetIndex = $i
    }
}

$anchorPara = $paras.Item($targetIndex)
$insPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml" />
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:b/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t>FRI Mar 08</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t xml:space="preserve"> 13:39:44 IST 2019</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t>Person Name</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
              <w:t>- CHANDRA</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t>---------------------------------------------------------------</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t>Item Name</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
              <w:t>- BEET</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Amount Received</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:tab/>
              <w:t>- 1078</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:t>Amount Received mode</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
              <w:tab/>
              <w:t>- CASH AND CLEARD</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PlainText"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insPoint.InsertXML($xml)
